$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 11, pushing the existing rows 11-89 down
# to 12-90 (preserving every original row's data/formatting as-is).
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44635
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101004
$ws.Range("J11").Value = "Frambuesa"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 690
$ws.Range("N11").Value = 7000
$ws.Range("O11").Value = 7500
$ws.Range("P11").Value = 7275
$ws.Range("Q11").Value = "$/bandeja 2 kilos"
$ws.Range("R11").Value = "Provincia de Linares"
$ws.Range("S11").Value = 3638
$ws.Range("T11").Value = 2
